# Append three new timesheet entries (rows 37-39) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Carry the date (col A) and duration (col B) formatting down from row 36 without
# dragging along its (wrap-text) comment formatting in col C.
$ws.Range("A36:B36").Copy()
$ws.Range("A37:B39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Excel serial date number for 2020-07-27 (row 36 was 44037 = 2020-07-25)
$newDateSerial = 44039

# Row 37: Konsultacje (existing shared string, reused text)
$ws.Cells.Item(37, 1).Value = $newDateSerial
$ws.Cells.Item(37, 2).Value = 0.5
$ws.Cells.Item(37, 3).Value = "Konsultacje"

# Row 38: new comment text
$ws.Cells.Item(38, 1).Value = $newDateSerial
$ws.Cells.Item(38, 2).Value = 1.5
$ws.Cells.Item(38, 3).Value = "internal_value dla XMLExportDC. Upgrade i instalacja VS2019"

# Row 39: new comment text
$ws.Cells.Item(39, 1).Value = $newDateSerial
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(39, 3).Value = "Instalacja ETCore. Budowa modelu istniejącej bazy. Tutoriale do EF."

# Update selection / active cell to mirror the saved view state
$ws.Range("C42").Select()

$wb.Save()
